# Update cell A3 on the "variables_####" worksheets that currently read
# "congenital" so that they read "misc_long_term" instead, matching the
# new dataset naming used by the baseline regression.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("A3")
    if ($cell.Value2 -eq "congenital") {
        $cell.Value2 = "misc_long_term"
    }
}
